# Revised templates with links to datasets in GBIF.org
#
# The README sheet gains two new lines (with a blank spacer line after
# them) right after the existing "reference" hyperlink line, pointing the
# reader to the GBIF.org page that indexes this example dataset. Every
# row below the insertion point shifts down by three.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("README")   # the "README" sheet (already the active tab)
$ws.Activate()

# Insert three blank rows right after row 8 (old rows 9-19 become 12-22).
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

# New row 9 reuses the plain-paragraph look of row 6 (A/B column styling,
# no explicit row height).
$ws.Range("A6:B6").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122) | Out-Null

# New rows 10 and 11 reuse the "reference link" look of row 8 (the blank
# row that currently trails the CUMV reference hyperlink in row 7).
$ws.Range("A8:B8").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:B11").PasteSpecial(-4122) | Out-Null

$ws.Range("B9").Value = "This dataset has been indexed by GBIF. To see how GBIF indexes the data go here:"
$ws.Range("B10").Value = "http://www.gbif.org/dataset/a8ee9bc6-5914-427d-9fba-f8545250ac34"
# B11 stays blank, mirroring B8.
